$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "244.52"
Set-TextCell "G2" "22"
Set-TextCell "G3" "22"
Set-TextCell "D4" "5.391"
Set-TextCell "G4" "22"
Set-TextCell "D5" "0.06048"
Set-TextCell "G5" "22"
Set-TextCell "D6" "3.396"
Set-TextCell "G6" "22"
Set-TextCell "D7" "0.8142"
Set-TextCell "G7" "22"
Set-TextCell "D8" "0.9228"
Set-TextCell "G8" "22"
Set-TextCell "B9" "WazirX"
Set-TextCell "C9" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell "D9" "0.1442"
Set-TextCell "E9" "8WazirXWRX"
Set-TextCell "G9" "22"
Set-TextCell "B10" "MandalaExchangeToken"
Set-TextCell "C10" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell "D10" "0.07450"
Set-TextCell "E10" "9MandalaExchangeTokenMDX"
Set-TextCell "G10" "22"
Set-TextCell "B11" "LiechtensteinCryptoassetsExchange"
Set-TextCell "C11" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell "D11" "0.03392"
Set-TextCell "E11" "10LiechtensteinCryptoassetsExchangeLCX"
Set-TextCell "G11" "22"
Set-TextCell "B12" "BitrueCoin"
Set-TextCell "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell "D12" "0.03049"
Set-TextCell "E12" "11BitrueCoinBTR"
Set-TextCell "G12" "22"
Set-TextCell "B13" "BitMartToken"
Set-TextCell "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell "D13" "0.09426"
Set-TextCell "E13" "12BitMartTokenBMX"
Set-TextCell "G13" "22"
Set-TextCell "B14" "MCDex"
Set-TextCell "C14" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell "D14" "4.010"
Set-TextCell "E14" "13MCDexMCB"
Set-TextCell "G14" "22"
Set-TextCell "B15" "BitForexToken"
Set-TextCell "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell "D15" "0.001591"
Set-TextCell "E15" "14BitForexTokenBF"
Set-TextCell "G15" "22"
Set-TextCell "B16" "CoinExToken"
Set-TextCell "C16" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell "D16" "0.04812"
Set-TextCell "E16" "15CoinExTokenCET"
Set-TextCell "G16" "22"
Set-TextCell "B17" "One"
Set-TextCell "C17" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell "D17" "0.0005946"
Set-TextCell "E17" "16OneONE"
Set-TextCell "G17" "22"
Set-TextCell "D18" "0.005415"
Set-TextCell "G18" "22"
Set-TextCell "G19" "22"
Set-TextCell "D20" "0.0009877"
Set-TextCell "G20" "22"
Set-TextCell "D21" "3.653"
Set-TextCell "G21" "22"
Set-TextCell "D22" "6.431"
Set-TextCell "G22" "22"
Set-TextCell "D23" "2.178"
Set-TextCell "G23" "22"
Set-TextCell "G24" "22"
Set-TextCell "G25" "22"
Set-TextCell "D26" "0.00008405"
Set-TextCell "G26" "22"
Set-TextCell "D27" "0.0002903"
Set-TextCell "G27" "22"
Set-TextCell "G28" "22"
Set-TextCell "G29" "22"
Set-TextCell "G30" "22"
Set-TextCell "G31" "22"
Set-TextCell "G32" "22"
Set-TextCell "G33" "22"
Set-TextCell "G34" "22"
Set-TextCell "G35" "22"
Set-TextCell "G36" "22"
Set-TextCell "G37" "22"
Set-TextCell "G38" "22"
Set-TextCell "G39" "22"
Set-TextCell "D40" "0.04022"
Set-TextCell "G40" "22"
Set-TextCell "B41" "BKEXToken"
Set-TextCell "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D41" "0.1075"
Set-TextCell "E41" "40BKEXTokenBKK"
Set-TextCell "G41" "22"
Set-TextCell "B42" "CEJI"
Set-TextCell "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D42" "0.002712"
Set-TextCell "E42" "41CEJICEJI"
Set-TextCell "G42" "22"
Set-TextCell "B43" "KickToken"
Set-TextCell "C43" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell "D43" "0.003048"
Set-TextCell "E43" "42KickTokenKICK"
Set-TextCell "G43" "22"
Set-TextCell "D44" "0.005833"
Set-TextCell "G44" "22"
Set-TextCell "G45" "22"
Set-TextCell "G46" "22"
Set-TextCell "E47" "46CoinbaseStockTokenCOINBestin24h"
Set-TextCell "G47" "22"
Set-TextCell "D48" "0.002321"
Set-TextCell "G48" "22"
Set-TextCell "D49" "0.00002102"
Set-TextCell "G49" "22"
Set-TextCell "G50" "22"
Set-TextCell "G51" "22"
